$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task "2. Make feature to clear old news" (row 3): status changes from "In Progress" to "Open"
$ws.Range("C3").Value = "Open"

# Task "17. Add links to sites" (row 18): status changes from "Open" to "Done" (with the green "Done" font color used elsewhere)
$ws.Range("C18").Value = "Done"
$ws.Range("C18").Font.Color = 5287936

# Update the active selection to reflect where the edit was made
$ws.Range("C14").Select()
